$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "material/geometry attributes of scatters." "material/geometry attributes of scatterers."
Replace-Text "constitutes the resolution of scattered wave fields" "constitutes the realization of scattered wave fields"
Replace-Text "Material constituents of the scatter give rise" "Material constituents of the scatterer give rise"
Replace-Text "no sources exist within scatters." "no sources exist within scatterers."
Replace-Text "matrix inversions, often using iterative Krylov Methods" "matrix inversions, or using iterative Krylov Methods"
Replace-Text "accelerate the matrix inversions by Fast Fourier Transforms" "accelerate the matrix multiplications by Fast Fourier Transforms"
Replace-Text "the material/geometry of the scatter as depicted" "the material/geometry of the scatterer as depicted"
Replace-Text "find the electric field given scatter geometry and material information" "find the electric field given scatterer geometry and material information"
